# Fruta / hortaliza, semanal
# Insert a new weekly record at row 481 (pushing the existing rows 481-530
# down to 482-531) on the active sheet of the "Uva" (grape) subset workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 481:530 down to 482:531, opening up a blank row 481.
$ws.Rows.Item(481).Insert()

# Populate the newly opened row 481 with the new observation.
$ws.Cells.Item(481, 1).Value  = 5
$ws.Cells.Item(481, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(481, 3).Value  = "Maule"
$ws.Cells.Item(481, 4).Value  = 44918
$ws.Cells.Item(481, 5).Value  = 7
$ws.Cells.Item(481, 6).Value  = "Fruta"
$ws.Cells.Item(481, 7).Value  = 100109
$ws.Cells.Item(481, 8).Value  = "Uva"
$ws.Cells.Item(481, 9).Value  = 100109001
$ws.Cells.Item(481, 10).Value = "Uva"
$ws.Cells.Item(481, 11).Value = "Superior Seedless"
$ws.Cells.Item(481, 12).Value = "Especial"
$ws.Cells.Item(481, 13).Value = 200
$ws.Cells.Item(481, 14).Value = 15000
$ws.Cells.Item(481, 15).Value = 15000
$ws.Cells.Item(481, 16).Value = 15000
$ws.Cells.Item(481, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(481, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(481, 19).Value = 1500
$ws.Cells.Item(481, 20).Value = 10
